$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B12").Value = 7.5
$ws.Range("B13").Value = 1071428.571428572
$ws.Range("B32").Value = 1171428.571428572
$ws.Range("B34").Value = 1171428.571428572
